# ref: only use http accessor to receive from API, not to send
#
# Updates the "Journal de travail" time-tracking sheet:
#  - bump a few existing "Temps [h]" entries
#  - add a new journal entry (row 71) for the API-key authentication work
#  - move the viewport/selection further down the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# --- Adjust existing hour entries -----------------------------------------
$ws.Range("C52").Value = 7      # was 6
$ws.Range("C55").Value = 6      # was 5.5
$ws.Range("C58").Value = 6      # was 5

# --- Add the new journal entry on row 71 -----------------------------------
$ws.Range("A71").Value = 45117
$ws.Range("B71").Value = "Implémentation"
$ws.Range("C71").Value = 9
$ws.Range("D71").Value = "Authentification par API keys pour l'envoi vers les modules"

# Recalculate the workbook so the SUM formulas in C77/D79 pick up the new values
$excel.Calculate()

# --- Update the view / selection state -------------------------------------
$win = $excel.ActiveWindow
$ws.Range("D81").Select()
$win.ScrollRow = 62
$win.ScrollColumn = 1
